$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue Tracking")

# Row 8 (Issue #6) - Resolution/Plan (E8): append the 01/18 status update
$ws.Range("E8").Value = @"
•Jimmy needs to have a teamviewer session to debug what is going on here.
•11/28 Teamviewer ready for Jimmy.
•11/29 Jimmy tested and saw a delay of 2 minutes when switching. He will investigate and see if there is any chance for improvement.
•There is a problem with this setup being offline every day. Jimmy is loosing time as he cannot test without an active TV.
•01/04 Jimmy is going to try to see the same issue in Belgium and if he does the TV setup will no longer be needed.
•01/11 The issue does not seem to be related to the firmware switch, but to something specific on the SIM. We are asking Gemalto to check this on the module.
•01/18 We are waiting feedback from Gemalto on this one. They are looking into it, but might take time.
"@

# Row 7 (Issue #5) - GetWireless comments (F7): new Brandon/SIM note
$ws.Range("F7").Value = @"
•01/13 Brandon saw a problem with the modem not reading the SIM correctly after switching it and rebooting the modem.
"@

# Row 7 (Issue #5) - Resolution/Plan (E7): append the 01/18 status update
$ws.Range("E7").Value = @"
•This is the same on both the current CG 3G (Gobi) and the CG LTE. There was no change compared to those devices.
•We will try to get to a solution for the release following the 2.70.0. Some testing will be needed from GetWireless and Option.
•Engineering build will be provided in January for testing and if OK, then a release will happen the same month.
•01/06 Engineering build provided to GetWireless. GetWireless to test.
•01/11 It will be great if we can get a result by the end of the week so that we can add it on the release candidate.
•01/18 Jimmy is looking into this one.
"@

# Row 8 (Issue #6) - OPTION INTERNAL COMMENTS (D8): append the 01/18 status update
$ws.Range("D8").Value = @"
•Jimmy requested a Teamviewer session. Franco to work on getting this setup.
•28/11 TV ready for Jimmy.
•29/11 Jimmy tested today and saw a delay of 2 minutes to find a signal. He will investigate.
•11/01 Does not seem to be related to the firmware switch. but to a specific SIM. We are asking Gemalto to check this on the modem.
•01/18 We are waiting feedback from Gemalto on this one. They are looking into it, but might take time.
"@

# Row 7 (Issue #5) - Status (G7): move from "Under Test (GW)" to "Under Investigation (OP)"
$ws.Range("G7").Value = "Under Investigation (OP)"

# Row 8 grew taller to fit the extra line of text
$ws.Rows.Item(8).RowHeight = 300

# Refresh the view: keep the header rows frozen and scroll down to row 8,
# leaving the active selection on E13
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("E13").Select()
